$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 5.553084769722144

# Row 3
$ws.Range("B3").Value = 0.127881588408715
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 1.742940831014585

# Row 4
$ws.Range("B4").Value = 0.6753301551942219
$ws.Range("C4").Value = 10.29869402782916
$ws.Range("D4").Value = 337.1190423067083
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 356.7532989756807
